$d = $word.ActiveDocument

# The document's final paragraph currently contains only the _GoBack
# bookmark (bookmarkStart immediately followed by bookmarkEnd, no runs).
# We split it into two paragraphs:
#   1) an empty paragraph (same "both" justification)
#   2) a paragraph that keeps the _GoBack bookmark and adds new text
#      describing the SICI20v2 backend folder / mongo backup info.
$lastPara = $d.Paragraphs.Last
$target = $lastPara.Range

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/></w:pPr></w:p>' +
          '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:pPr><w:jc w:val="both"/></w:pPr>' +
          '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
          '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">En el </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:rPr><w:b/></w:rPr><w:t>backend</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve"> hay una carpeta llamada SICI20v2</w:t></w:r>' +
          '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
          '<w:bookmarkEnd w:id="0"/>' +
          '<w:r><w:t xml:space="preserve">ese es el </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>backup</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t xml:space="preserve"> de la base de datos mongo para que hablemos el mismo idioma, todos los datos que contiene son </w:t></w:r>' +
          '<w:proofErr w:type="spellStart"/>' +
          '<w:r><w:t>dtos</w:t></w:r>' +
          '<w:proofErr w:type="spellEnd"/>' +
          '<w:r><w:t xml:space="preserve"> de prueba.</w:t></w:r>' +
          '</w:p>'

$null = $target.InsertXML($newXml)
